$wb = $excel.ActiveWorkbook

# --- Worksheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2600
$ws.Range("I40").Value = 2400
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2400
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2225
$ws.Range("N40").Value = -3350
# Row 96
$ws.Range("H96").Value = 1360.9048
$ws.Range("I96").Value = 506.69232
$ws.Range("K96").Value = 1520.07696
$ws.Range("M96").Value = -147.0769599999999
# Row 112
$ws.Range("H112").Value = 1968.9333
$ws.Range("I112").Value = 1673.25
$ws.Range("J112").Value = 2076.4546
$ws.Range("K112").Value = 5019.75
$ws.Range("L112").Value = 6229.3638
$ws.Range("M112").Value = -3911.75
$ws.Range("N112").Value = -8445.363799999999
# Row 127
$ws.Range("H127").Value = 337000
$ws.Range("I127").Value = 337000
$ws.Range("K127").Value = 1011000
$ws.Range("M127").Value = -1006040
# Row 129
$ws.Range("H129").Value = 1193.4445
$ws.Range("I129").Value = 863
$ws.Range("J129").Value = 2350
$ws.Range("K129").Value = 2589
$ws.Range("L129").Value = 7050
$ws.Range("M129").Value = 2411
$ws.Range("N129").Value = -17050
# Row 131
$ws.Range("H131").Value = 2460
$ws.Range("I131").Value = 2460
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 7380
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -2340
$ws.Range("N131").ClearContents()
# Row 137
$ws.Range("H137").Value = 13527282
$ws.Range("I137").Value = 27780860
$ws.Range("J137").Value = 23892.422
$ws.Range("K137").Value = 83342580
$ws.Range("L137").Value = 71677.266
$ws.Range("M137").Value = -83340030
$ws.Range("N137").Value = -76777.266
# Row 138
$ws.Range("H138").Value = 7941.838
$ws.Range("I138").Value = 7760.8887
$ws.Range("K138").Value = 23282.6661
$ws.Range("M138").Value = -18142.6661
# Row 141
$ws.Range("H141").Value = 11915.346
$ws.Range("I141").Value = 3699.8333
$ws.Range("J141").Value = 14380
$ws.Range("K141").Value = 11099.4999
$ws.Range("L141").Value = 43140
$ws.Range("M141").Value = -5919.499899999999
$ws.Range("N141").Value = -53500

# --- Worksheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 583.625
$ws.Range("I5").Value = 123.8
$ws.Range("J5").Value = 1350
$ws.Range("K5").Value = 123.8
$ws.Range("L5").Value = 1350
$ws.Range("M5").Value = -11.8
$ws.Range("N5").Value = -1574
# Row 45
$ws.Range("H45").Value = 1749.4
$ws.Range("I45").Value = 1377.1111
$ws.Range("J45").Value = 5100
$ws.Range("K45").Value = 1377.1111
$ws.Range("L45").Value = 5100
$ws.Range("M45").Value = -1000.1111
$ws.Range("N45").Value = -5854

# --- Worksheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 583.625
$ws.Range("I4").Value = 123.8
$ws.Range("J4").Value = 1350
$ws.Range("K4").Value = 123.8
$ws.Range("L4").Value = 1350
$ws.Range("M4").Value = -8.799999999999997
$ws.Range("N4").Value = -1580
# Row 8
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
# Row 14
$ws.Range("H14").Value = 28000
$ws.Range("J14").Value = 28000
$ws.Range("L14").Value = 28000
$ws.Range("N14").Value = -28344

# --- Worksheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 566416.5600000001
$ws.Range("I31").Value = 797995.5600000001
$ws.Range("J31").Value = 4010.5
$ws.Range("K31").Value = 797995.5600000001
$ws.Range("L31").Value = 4010.5
$ws.Range("M31").Value = -797700.5600000001
$ws.Range("N31").Value = -4600.5
# Row 34
$ws.Range("H34").Value = 566416.5600000001
$ws.Range("I34").Value = 797995.5600000001
$ws.Range("J34").Value = 4010.5
$ws.Range("K34").Value = 797995.5600000001
$ws.Range("L34").Value = 4010.5
$ws.Range("M34").Value = -797793.5600000001
$ws.Range("N34").Value = -4414.5
# Row 107
$ws.Range("H107").Value = 461.35294
$ws.Range("I107").Value = 458.9375
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 458.9375
$ws.Range("L107").Value = 500
$ws.Range("M107").Value = 1461.0625
$ws.Range("N107").Value = -4340
# Row 132
$ws.Range("H132").Value = 2285
$ws.Range("I132").Value = 2227.3572
$ws.Range("J132").Value = 3899
$ws.Range("K132").Value = 6682.071599999999
$ws.Range("L132").Value = 11697
$ws.Range("M132").Value = -4152.071599999999
$ws.Range("N132").Value = -16757

# --- Worksheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 339.5625
$ws.Range("I12").Value = 1262
$ws.Range("J12").Value = 32.083332
$ws.Range("K12").Value = 3786
$ws.Range("L12").Value = 96.249996
$ws.Range("M12").Value = -3613
$ws.Range("N12").Value = -442.249996
# Row 107
$ws.Range("H107").Value = 4554.154
$ws.Range("J107").Value = 6100.5
$ws.Range("L107").Value = 18301.5
$ws.Range("N107").Value = -22141.5
# Row 129
$ws.Range("H129").Value = 2634.4
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2634.4
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 7903.200000000001
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -17903.2
# Row 131
$ws.Range("H131").Value = 4182.354
$ws.Range("J131").Value = 5030.3784
$ws.Range("L131").Value = 15091.1352
$ws.Range("N131").Value = -25171.1352
# Row 137
$ws.Range("H137").Value = 8316.333000000001
$ws.Range("I137").Value = 2006.4286
$ws.Range("J137").Value = 13837.5
$ws.Range("K137").Value = 6019.2858
$ws.Range("L137").Value = 41512.5
$ws.Range("M137").Value = -919.2857999999997
$ws.Range("N137").Value = -51712.5

# --- Worksheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1008.69446
$ws.Range("I97").Value = 826.5833
$ws.Range("J97").Value = 1372.9166
$ws.Range("K97").Value = 826.5833
$ws.Range("L97").Value = 1372.9166
$ws.Range("M97").Value = -330.5833
$ws.Range("N97").Value = -2364.9166

# --- Worksheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 318.9737
$ws.Range("I55").Value = 119.818184
$ws.Range("J55").Value = 592.8125
$ws.Range("K55").Value = 119.818184
$ws.Range("L55").Value = 592.8125
$ws.Range("M55").Value = 53.181816
$ws.Range("N55").Value = -938.8125
# Row 93
$ws.Range("H93").Value = 2793.28
$ws.Range("I93").Value = 1645.0834
$ws.Range("J93").Value = 3853.1538
$ws.Range("K93").Value = 1645.0834
$ws.Range("L93").Value = 3853.1538
$ws.Range("M93").Value = -397.0834
$ws.Range("N93").Value = -6349.1538
# Row 100
$ws.Range("H100").Value = 2507.5789
$ws.Range("I100").Value = 2209.6667
$ws.Range("J100").Value = 3624.75
$ws.Range("K100").Value = 2209.6667
$ws.Range("L100").Value = 3624.75
$ws.Range("M100").Value = -1668.6667
$ws.Range("N100").Value = -4706.75
